# Generate Report for Handoff
# Updates the localization status report: cells previously showing
# "In Translation" move to "Ready for handoff", and the associated
# timestamp cells are bumped to reflect the new handoff generation time.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status columns: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handback DateTime (de-de) bumped together
$wsOverview.Range("G2").Value = "2016-11-15 16:16:19"
$wsDeDe.Range("H2").Value     = "2016-11-15 16:16:19"

# Latest Handoff Datetime (zh-cn) bumped
$wsZhCn.Range("H2").Value = "2016-11-15 16:16:05"

# Auto-fit the status columns now that they hold the longer text
$wsOverview.Range("E1:F2").Columns.AutoFit() | Out-Null
$wsZhCn.Range("C1:C2").Columns.AutoFit() | Out-Null
$wsDeDe.Range("C1:C2").Columns.AutoFit() | Out-Null
